# Add a new "Init" worksheet after the existing "GoHome" sheet, mirroring the
# same 8-column x 22-row template layout used by "Login"/"GoHome", and fill in
# the single "Navigator / Open" action row (row 3).

$wb = $excel.ActiveWorkbook

# --- Create the new sheet at the end of the workbook -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Init"

# --- Column widths (matches the Login/GoHome template) ---------------------
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 8
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 12
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 10

# --- Stamp every cell of the A1:H22 template grid with a style so the sheet
# keeps the same "fully-formatted" 22x8 shape as Login/GoHome -------------
for ($r = 1; $r -le 22; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Borders.Item(7).Weight = 2
    }
}

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "Flow"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Object"
$ws.Range("D1").Value = "Action"
$ws.Range("E1").Value = "Param Name"
$ws.Range("F1").Value = "Param Type"
$ws.Range("G1").Value = "Param Value"

# --- Row 3: the single "Open" navigator action --------------------------
$ws.Range("B3").Value = "Action"
$ws.Range("C3").Value = "Navigator"
$ws.Range("D3").Value = "Open"
$ws.Range("E3").Value = "url"
$ws.Range("F3").Value = "string"
$ws.Range("G3").Value = "https://inflectra365.crm.dynamics.com/main.aspx?appid=b703cc78-b50e-ea11-a812-000d3a8c9a6d"

# --- Restore the originally-active sheet/tab ----------------------------
$wb.Worksheets.Item("Login").Activate()
